# "separated Start for testing, added report message to Mass_Balance_Cor"
#
# Rescales the six flow-rate readings in column E (rows 20-25) of Sheet1
# from raw integer units down to liters (divide by 1000) and leaves the
# live selection on the cell the reviewer ended up at (G25) while testing,
# instead of the old J10:J12 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rescale E20:E25 (divide the raw integer reading by 1000) ----------
# NB: read back via Value2 - Value's getter is unreliable through this
# COM shim (it hands back the member descriptor instead of the scalar).
$rows = 20..25
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $raw = $cell.Value2
    $cell.Value = $raw / 1000
}

# --- Move the selection to where review/testing continued --------------
$ws.Activate()
[void]$ws.Range("G25").Select()

Write-Output "Mass_Balance_Cor rescaled E20:E25 and moved selection to G25"
